# Add a "password" column (D) to the user-list sheet, populating the
# header (bold, like the other headers) and the per-row password values
# that line up with the existing username/nama rows:
#   row2 -> admin_alex     -> alexadmin1234567
#   row3 -> KucingImut     -> kucing1234567
#   row4 -> lele           -> lele1234567

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell: bold font, same shared-string/style family as A1:C1.
$ws.Range("D1").Value = "password"
$ws.Range("D1").Font.Bold = $true

# Data cells - written bottom-up so the shared-strings table ends up in
# the same order (password, lele1234567, kucing1234567, alexadmin1234567)
# as the authored workbook.
$ws.Range("D4").Value = "lele1234567"
$ws.Range("D3").Value = "kucing1234567"
$ws.Range("D2").Value = "alexadmin1234567"

# New column gets its own (slightly narrower) width, matching the other
# data columns' style of an explicit custom width.
$ws.Columns.Item(4).ColumnWidth = 13.833333333333332

# Leave the selection on the newly added cell, as in the authored file.
$ws.Range("D3").Select()
